$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new data rows above the last (bottom-bordered) data row ---
# Old row 19 ("DEIMER...") shifts down to row 22, and the signature block
# (old rows 24/25) shifts down to rows 27/28, carrying its own formatting
# with it automatically.
$ws.Range("B19:J21").Insert()

# Copy the formatting (borders, fonts, number formats) of the regular data
# row (row 18) onto the three freshly inserted rows so they match the rest
# of the table instead of picking up the default/no-border style.
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update summary header values ---
$ws.Range("E11").Value = 203290
$ws.Range("C13").Value = 6
$ws.Range("F13").Value = 6

# --- Row 16: FRANCISCO JOSE GOMEZ PAJARO ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "9149166"
$ws.Range("D16").Value = "FRANCISCO JOSE GOMEZ PAJARO"
$ws.Range("E16").Value = "1807"
$ws.Range("F16").Value = 31249
$ws.Range("G16").Value = 801000

# --- Row 17: JOSE BARRIOS SANCHEZ ---
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73169213"
$ws.Range("D17").Value = "JOSE BARRIOS SANCHEZ"
$ws.Range("E17").Value = "2102"
$ws.Range("F17").Value = 41400
$ws.Range("G17").Value = 1035000

# --- Row 18: JOSE ENRIQUE ALVAREZ ESCOBAR ---
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "9145804"
$ws.Range("D18").Value = "JOSE ENRIQUE ALVAREZ ESCOBAR"
$ws.Range("E18").Value = "2104"
$ws.Range("F18").Value = 33120
$ws.Range("G18").Value = 1100000

# --- Row 19: JAIRO RIOS MENDOZA (periodo 1807) ---
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1047462719"
$ws.Range("D19").Value = "JAIRO RIOS MENDOZA"
$ws.Range("E19").Value = "1807"
$ws.Range("F19").Value = 31249
$ws.Range("G19").Value = 781242

# --- Row 20: JAIRO RIOS MENDOZA (periodo 1806) ---
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1047462719"
$ws.Range("D20").Value = "JAIRO RIOS MENDOZA"
$ws.Range("E20").Value = "1806"
$ws.Range("F20").Value = 31249
$ws.Range("G20").Value = 781242

# --- Row 21: JHON ALFREDO JUNCO CORREA ---
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1193515790"
$ws.Range("D21").Value = "JHON ALFREDO JUNCO CORREA"
$ws.Range("E21").Value = "2507"
$ws.Range("F21").Value = 1898
$ws.Range("G21").Value = 1423500

# --- Row 22: DEIMER ANTONIO MARRUGO HERRERA (unchanged values, already shifted here) ---
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1047470076"
$ws.Range("D22").Value = "DEIMER ANTONIO MARRUGO HERRERA"
$ws.Range("E22").Value = "1908"
$ws.Range("F22").Value = 33125
$ws.Range("G22").Value = 828116
